# Regenerate the localization-status report for a new handoff:
#  - row 2 file "aeed826a-d438-4995-8efd-3e07451c46f9.md" -> "ec347a14-2d49-483f-b70a-b82f66d3c54f.md"
#  - row 3 file "d65be23c-4c3c-4f2e-b6f5-8d98b3b3b70e.md" -> "fffff0f53630-fc9e-4919-b915-65fb2cb8d0e4.md"
#  - status "Handed back: in sync with en-US" -> "Ready for handoff"
#  - new xliff hashes / timestamps, handback info cleared (not yet handed back)

$wb = $excel.ActiveWorkbook

$oldGuid1 = "aeed826a-d438-4995-8efd-3e07451c46f9"
$newGuid1 = "ec347a14-2d49-483f-b70a-b82f66d3c54f"
$oldGuid2 = "d65be23c-4c3c-4f2e-b6f5-8d98b3b3b70e"
$newGuid2 = "fffff0f53630-fc9e-4919-b915-65fb2cb8d0e4"
$newHash  = "5b45b90031a49edabfb6b20291b1375b956b5d29"

$repoBase = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ca3db49d09671fd1930707dee7f134de8063a4f6/e2e"

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

$ov.Range("A2").Value = "$newGuid1.md"
$ov.Range("E2").Value = "Ready for handoff"
$ov.Range("F2").Value = "Ready for handoff"
$ov.Range("G2").Value = "2016-08-29 13:03:36"

$ov.Range("A3").Value = "$newGuid2.md"
$ov.Range("E3").Value = "Ready for handoff"
$ov.Range("F3").Value = "Ready for handoff"
$ov.Range("G3").Value = "2016-08-29 13:03:36"

$ov.Hyperlinks.Delete()
$ov.Hyperlinks.Add($ov.Range("B2"), "$repoBase/$newGuid1.md", [Type]::Missing, [Type]::Missing, "e2e\$newGuid1.md")
$ov.Hyperlinks.Add($ov.Range("B3"), "$repoBase/$newGuid2.md", [Type]::Missing, [Type]::Missing, "e2e\$newGuid2.md")

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("A2").Value = "$newGuid1.md"
$zh.Range("C2").Value = "Ready for handoff"
$zh.Range("G2").Value = "$newGuid1.$newHash.zh-cn.xlf"
$zh.Range("H2").Value = "2016-08-29 13:03:31"
$zh.Range("I2").Value = ""
$zh.Range("I2").Style = "Normal"
$zh.Range("J2").Value = ""
$zh.Range("K2").Value = "0001-01-01 00:00:00"

$zh.Range("A3").Value = "$newGuid2.md"
$zh.Range("C3").Value = "Ready for handoff"
$zh.Range("F3").Value = "'True"
$zh.Range("F3").Style = "Normal"
$zh.Range("G3").Value = "$newGuid1.$newHash.zh-cn.xlf"
$zh.Range("H3").Value = "2016-08-29 13:03:31"
$zh.Range("I3").Value = ""
$zh.Range("I3").Style = "Normal"
$zh.Range("J3").Value = ""
$zh.Range("K3").Value = "0001-01-01 00:00:00"

$zh.Hyperlinks.Delete()
$zh.Hyperlinks.Add($zh.Range("A2"), "$repoBase/$newGuid1.md", [Type]::Missing, [Type]::Missing, "$newGuid1.md")
$zh.Hyperlinks.Add($zh.Range("A3"), "$repoBase/$newGuid2.md", [Type]::Missing, [Type]::Missing, "$newGuid2.md")

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("A2").Value = "$newGuid1.md"
$de.Range("C2").Value = "Ready for handoff"
$de.Range("G2").Value = "$newGuid1.$newHash.de-de.xlf"
$de.Range("H2").Value = "2016-08-29 13:03:36"
$de.Range("I2").Value = ""
$de.Range("I2").Style = "Normal"
$de.Range("J2").Value = ""
$de.Range("K2").Value = "0001-01-01 00:00:00"

$de.Range("A3").Value = "$newGuid2.md"
$de.Range("C3").Value = "Ready for handoff"
$de.Range("F3").Value = "'True"
$de.Range("F3").Style = "Normal"
$de.Range("G3").Value = "$newGuid1.$newHash.de-de.xlf"
$de.Range("H3").Value = "2016-08-29 13:03:36"
$de.Range("I3").Value = ""
$de.Range("I3").Style = "Normal"
$de.Range("J3").Value = ""
$de.Range("K3").Value = "0001-01-01 00:00:00"

$de.Hyperlinks.Delete()
$de.Hyperlinks.Add($de.Range("A2"), "$repoBase/$newGuid1.md", [Type]::Missing, [Type]::Missing, "$newGuid1.md")
$de.Hyperlinks.Add($de.Range("A3"), "$repoBase/$newGuid2.md", [Type]::Missing, [Type]::Missing, "$newGuid2.md")

# ---------------------------------------------------------------------------
# Column widths: auto-fit now that some long filenames collapsed to blank
# ---------------------------------------------------------------------------
$ov.Columns.Item(5).AutoFit() | Out-Null
$ov.Columns.Item(6).AutoFit() | Out-Null
$zh.Columns.Item(3).AutoFit() | Out-Null
$zh.Columns.Item(9).AutoFit() | Out-Null
$zh.Columns.Item(10).AutoFit() | Out-Null
$de.Columns.Item(3).AutoFit() | Out-Null
$de.Columns.Item(9).AutoFit() | Out-Null
$de.Columns.Item(10).AutoFit() | Out-Null
